$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 818.375
$ws.Range("I19").Value = 800
$ws.Range("J19").Value = 829.4
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 829.4
$ws.Range("M19").Value = -625
$ws.Range("N19").Value = -1179.4
$ws.Range("H38").Value = 1927.1111
$ws.Range("J38").Value = 3999
$ws.Range("L38").Value = 11997
$ws.Range("N38").Value = -12741
$ws.Range("H80").Value = 400.77777
$ws.Range("I80").Value = 172.25
$ws.Range("J80").Value = 583.6
$ws.Range("K80").Value = 516.75
$ws.Range("L80").Value = 1750.8
$ws.Range("M80").Value = 481.25
$ws.Range("N80").Value = -3746.8
$ws.Range("H83").Value = 400.77777
$ws.Range("I83").Value = 172.25
$ws.Range("J83").Value = 583.6
$ws.Range("K83").Value = 1550.25
$ws.Range("L83").Value = 5252.400000000001
$ws.Range("M83").Value = 3441.75
$ws.Range("N83").Value = -15236.4
$ws.Range("H92").Value = 2259.2
$ws.Range("I92").Value = 2560.5
$ws.Range("J92").Value = 1749.3077
$ws.Range("K92").Value = 2560.5
$ws.Range("L92").Value = 1749.3077
$ws.Range("M92").Value = -1312.5
$ws.Range("N92").Value = -4245.3077
$ws.Range("H94").Value = 967.625
$ws.Range("I94").Value = 1183.1428
$ws.Range("K94").Value = 1183.1428
$ws.Range("M94").Value = -732.1428000000001
$ws.Range("H100").Value = 3105.0833
$ws.Range("I100").Value = 2862.6
$ws.Range("J100").Value = 3509.2222
$ws.Range("K100").Value = 2862.6
$ws.Range("L100").Value = 3509.2222
$ws.Range("M100").Value = -2321.6
$ws.Range("N100").Value = -4591.2222
$ws.Range("H112").Value = 6138.0938
$ws.Range("I112").Value = 70000
$ws.Range("J112").Value = 4078.0322
$ws.Range("K112").Value = 210000
$ws.Range("L112").Value = 12234.0966
$ws.Range("M112").Value = -208892
$ws.Range("N112").Value = -14450.0966
$ws.Range("H138").Value = 2328.9355
$ws.Range("I138").Value = 1761.75
$ws.Range("J138").Value = 2526.2173
$ws.Range("K138").Value = 5285.25
$ws.Range("L138").Value = 7578.651899999999
$ws.Range("M138").Value = -145.25
$ws.Range("N138").Value = -17858.6519

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2710.7144
$ws.Range("I32").Value = 1888.3137
$ws.Range("K32").Value = 1888.3137
$ws.Range("M32").Value = -1601.3137
$ws.Range("H44").Value = 48949
$ws.Range("J44").Value = 48949
$ws.Range("L44").Value = 48949
$ws.Range("N44").Value = -49925
$ws.Range("H45").Value = 3157.1667
$ws.Range("I45").Value = 3098.1177
$ws.Range("J45").Value = 3300.5715
$ws.Range("K45").Value = 3098.1177
$ws.Range("L45").Value = 3300.5715
$ws.Range("M45").Value = -2721.1177
$ws.Range("N45").Value = -4054.5715
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = ""
$ws.Range("N54").Value = 0
$ws.Range("H62").Value = 79124.5
$ws.Range("J62").Value = 88249
$ws.Range("L62").Value = 88249
$ws.Range("N62").Value = -89497
$ws.Range("H65").Value = 79124.5
$ws.Range("J65").Value = 88249
$ws.Range("L65").Value = 264747
$ws.Range("N65").Value = -270987
$ws.Range("H88").Value = 1694
$ws.Range("I88").Value = 1569
$ws.Range("J88").Value = 1746.0834
$ws.Range("K88").Value = 1569
$ws.Range("L88").Value = 1746.0834
$ws.Range("M88").Value = -1163
$ws.Range("N88").Value = -2558.0834
$ws.Range("H91").Value = 1694
$ws.Range("I91").Value = 1569
$ws.Range("J91").Value = 1746.0834
$ws.Range("K91").Value = 1569
$ws.Range("L91").Value = 1746.0834
$ws.Range("M91").Value = -165
$ws.Range("N91").Value = -4554.0834
$ws.Range("H110").Value = 4550882.5
$ws.Range("I110").Value = 4550882.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 4550882.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = -4548837.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6076405
$ws.Range("I20").Value = 11115341
$ws.Range("K20").Value = 11115341
$ws.Range("M20").Value = -11115094
$ws.Range("H105").Value = 111112744
$ws.Range("I105").Value = 111112744
$ws.Range("K105").Value = 111112744
$ws.Range("M105").Value = -111110997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2950.4
$ws.Range("I19").Value = 2950.4
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 2950.4
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -2780.4
$ws.Range("H24").Value = 2950.4
$ws.Range("I24").Value = 2950.4
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 2950.4
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = -2780.4
$ws.Range("H93").Value = 10353.5
$ws.Range("I93").Value = 10353.5
$ws.Range("K93").Value = 10353.5
$ws.Range("M93").Value = -8481.5
$ws.Range("H105").Value = 8173.3335
$ws.Range("J105").Value = 1329.5
$ws.Range("L105").Value = 1329.5
$ws.Range("N105").Value = -4823.5
$ws.Range("H122").Value = 1759.5186
$ws.Range("I122").Value = 1478.2609
$ws.Range("J122").Value = 3376.75
$ws.Range("K122").Value = 4434.7827
$ws.Range("L122").Value = 10130.25
$ws.Range("M122").Value = -1984.7827
$ws.Range("N122").Value = -15030.25
$ws.Range("H134").Value = 26321070
$ws.Range("I134").Value = 1697.7916
$ws.Range("K134").Value = 5093.3748
$ws.Range("M134").Value = -2558.3748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.692307
$ws.Range("I2").Value = 45.23077
$ws.Range("J2").Value = 58.153847
$ws.Range("K2").Value = 271.38462
$ws.Range("L2").Value = 348.923082
$ws.Range("M2").Value = -158.38462
$ws.Range("N2").Value = -574.923082
$ws.Range("H62").Value = 3607.6365
$ws.Range("I62").Value = 1585.5
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 4756.5
$ws.Range("L62").Value = 27000
$ws.Range("M62").Value = -4070.5
$ws.Range("N62").Value = -28372
$ws.Range("H65").Value = 3607.6365
$ws.Range("I65").Value = 1585.5
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 14269.5
$ws.Range("L65").Value = 81000
$ws.Range("M65").Value = -10837.5
$ws.Range("N65").Value = -87864
$ws.Range("H117").Value = 826.7143
$ws.Range("J117").Value = 1399.5
$ws.Range("L117").Value = 4198.5
$ws.Range("N117").Value = -11082.5
$ws.Range("H129").Value = 8265391
$ws.Range("I129").Value = 1030.2858
$ws.Range("K129").Value = 3090.8574
$ws.Range("M129").Value = 1909.1426
$ws.Range("H137").Value = 3399.6667
$ws.Range("I137").Value = 2402.5715
$ws.Range("J137").Value = 4795.6
$ws.Range("K137").Value = 7207.7145
$ws.Range("L137").Value = 14386.8
$ws.Range("M137").Value = -2107.7145
$ws.Range("N137").Value = -24586.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 16384.908
$ws.Range("J58").Value = 15137.667
$ws.Range("L58").Value = 15137.667
$ws.Range("N58").Value = -15691.667
$ws.Range("H70").Value = 14729.214
$ws.Range("I70").Value = 17150.625
$ws.Range("K70").Value = 17150.625
$ws.Range("M70").Value = -16880.625
$ws.Range("H73").Value = 14729.214
$ws.Range("I73").Value = 17150.625
$ws.Range("K73").Value = 17150.625
$ws.Range("M73").Value = -16214.625
$ws.Range("H97").Value = 869.8421
$ws.Range("I97").Value = 850.24243
$ws.Range("K97").Value = 850.24243
$ws.Range("M97").Value = -354.24243
$ws.Range("H122").Value = 1416139.6
$ws.Range("I122").Value = 1788229
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5364687
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -5362237
$ws.Range("N122").Value = -11500
$ws.Range("H126").Value = 3840108.5
$ws.Range("I126").Value = 2401595
$ws.Range("K126").Value = 7204785
$ws.Range("M126").Value = -7202315

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 83336280
$ws.Range("I16").Value = 100003360
$ws.Range("J16").Value = 871
$ws.Range("K16").Value = 100003360
$ws.Range("L16").Value = 871
$ws.Range("M16").Value = -100003190
$ws.Range("N16").Value = -1211
$ws.Range("H70").Value = 70000
$ws.Range("J70").Value = 70000
$ws.Range("L70").Value = 70000
$ws.Range("N70").Value = -70540
$ws.Range("H73").Value = 70000
$ws.Range("J73").Value = 70000
$ws.Range("L73").Value = 70000
$ws.Range("N73").Value = -71872
$ws.Range("H100").Value = 5358
$ws.Range("J100").Value = 2997
$ws.Range("L100").Value = 2997
$ws.Range("N100").Value = -4079
$ws.Range("H122").Value = 20994332
$ws.Range("I122").Value = 37199124
$ws.Range("J122").Value = 2088737.9
$ws.Range("K122").Value = 111597372
$ws.Range("L122").Value = 6266213.699999999
$ws.Range("M122").Value = -111594922
$ws.Range("N122").Value = -6271113.699999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5610601
$ws.Range("I11").Value = 9337001
$ws.Range("J11").Value = 21000
$ws.Range("K11").Value = 9337001
$ws.Range("L11").Value = 21000
$ws.Range("M11").Value = -9336859
$ws.Range("N11").Value = -21284
$ws.Range("H122").Value = 329835.44
$ws.Range("I122").Value = 421468.44
$ws.Range("K122").Value = 1264405.32
$ws.Range("M122").Value = -1261955.32
